$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.62"
$ws.Range("E2").Value = "'1.86%"
$ws.Range("D3").Value = "'27.20"
$ws.Range("E3").Value = "'1.48%"
$ws.Range("D4").Value = "'4.907"
$ws.Range("E5").Value = "'1.47%"
$ws.Range("E6").Value = "'0.61%"
$ws.Range("D7").Value = "'1.244"
$ws.Range("E7").Value = "'-4.15%"
$ws.Range("D8").Value = "'0.8839"
$ws.Range("E8").Value = "'-0.37%"
$ws.Range("E9").Value = "'4.16%"
$ws.Range("D10").Value = "'0.05011"
$ws.Range("E10").Value = "'-2.62%"
$ws.Range("D11").Value = "'0.07516"
$ws.Range("E11").Value = "'1.73%"
$ws.Range("E12").Value = "'-8.00%"
$ws.Range("D13").Value = "'0.09001"
$ws.Range("E13").Value = "'-0.40%"
$ws.Range("D14").Value = "'0.001575"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("D15").Value = "'0.0006413"
$ws.Range("E15").Value = "'1.82%"
$ws.Range("D16").Value = "'0.005818"
$ws.Range("E16").Value = "'-3.66%"
$ws.Range("D17").Value = "'3.459"
$ws.Range("E17").Value = "'-0.04%"
$ws.Range("D18").Value = "'3.314"
$ws.Range("E18").Value = "'-1.33%"
$ws.Range("D20").Value = "'0.3136"
$ws.Range("E20").Value = "'-0.94%"
$ws.Range("D21").Value = "'0.1336"
$ws.Range("E21").Value = "'0.07%"
$ws.Range("D22").Value = "'3.905"
$ws.Range("E22").Value = "'-0.18%"
$ws.Range("D23").Value = "'0.04439"
$ws.Range("E23").Value = "'1.79%"
$ws.Range("D24").Value = "'0.001173"
$ws.Range("E24").Value = "'-0.30%"
$ws.Range("E25").Value = "'5.18%"
$ws.Range("E26").Value = "'-0.12%"
$ws.Range("E27").Value = "'13.91%"
$ws.Range("D40").Value = "'0.04142"
$ws.Range("E40").Value = "'2.92%"
$ws.Range("D41").Value = "'0.006798"
$ws.Range("E41").Value = "'2.80%"
$ws.Range("D42").Value = "'0.1177"
$ws.Range("E42").Value = "'1.19%"
$ws.Range("E43").Value = "'13.71%"
$ws.Range("E44").Value = "'-3.96%"
$ws.Range("D45").Value = "'0.00005206"
$ws.Range("E45").Value = "'-2.21%"
$ws.Range("D46").Value = "'1.487"
$ws.Range("E46").Value = "'-36.96%"
$ws.Range("D47").Value = "'0.02024"
$ws.Range("E47").Value = "'-22.24%"
